$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# --- Update the "01/08/2023" style dates (column B, dd/mm/yyyy text) to "03/08/2023" ---
# These cells already use the quote-prefixed Text number format (@), so plain
# re-assignment (with a leading apostrophe to keep them as literal text) is enough.
$dmyCells = @("B1","B2","B3","B4","B5","B6","B7","B8","B10","B11")
foreach ($ref in $dmyCells) {
    $ws.Range($ref).Value = "'03/08/2023"
}

# B16 / B17 currently carry the quote-prefixed Date number format (like B12/B13/B14);
# the edit also restyles them to the quote-prefixed Text format used by the rest of
# the column, so force the number format before writing the new text.
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "'03/08/2023"

$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "'03/08/2023"

# --- Update the "2023/08/01" (yyyy/mm/dd) cell ---
$ws.Range("B12").Value = "'2023/08/03"

# --- Update the "08/01/2023" (mm/dd/yyyy) cell ---
$ws.Range("B14").Value = "'08/03/2023"

# --- Move the selection / scroll position (horizontal navigation bar scroll) ---
$ws.Range("B17").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1 | Out-Null
